# Update "想去人数" (want-to-go count) and one "最低票价" (min ticket price)
# figures across the 展览 / 演出 / 全部类型 sheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1856
$ws1.Range("F5").Value  = 10
$ws1.Range("F6").Value  = 817
$ws1.Range("G6").Value  = 50
$ws1.Range("F9").Value  = 38
$ws1.Range("F10").Value = 107
$ws1.Range("F16").Value = 4366
$ws1.Range("F19").Value = 477
$ws1.Range("F21").Value = 8
$ws1.Range("F22").Value = 993
$ws1.Range("F23").Value = 1768
$ws1.Range("F24").Value = 367
$ws1.Range("F26").Value = 16
$ws1.Range("F27").Value = 47
$ws1.Range("F28").Value = 2057
$ws1.Range("F29").Value = 71
$ws1.Range("F30").Value = 64
$ws1.Range("F32").Value = 145
$ws1.Range("F33").Value = 93
$ws1.Range("F34").Value = 24
$ws1.Range("F35").Value = 212
$ws1.Range("F36").Value = 29

# --- Sheet "演出" (sheetId 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 10

# --- Sheet "全部类型" (sheetId 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1856
$ws4.Range("F5").Value  = 10
$ws4.Range("F6").Value  = 817
$ws4.Range("F9").Value  = 38
$ws4.Range("F10").Value = 107
$ws4.Range("F16").Value = 10
$ws4.Range("F17").Value = 4366
$ws4.Range("F20").Value = 477
$ws4.Range("F22").Value = 8
$ws4.Range("F23").Value = 993
$ws4.Range("F24").Value = 1768
$ws4.Range("F25").Value = 367
$ws4.Range("F27").Value = 16
$ws4.Range("F28").Value = 47
$ws4.Range("F29").Value = 2057
$ws4.Range("F30").Value = 71
$ws4.Range("F31").Value = 64
$ws4.Range("F33").Value = 145
$ws4.Range("F34").Value = 93
$ws4.Range("F35").Value = 24
$ws4.Range("F36").Value = 212
$ws4.Range("F37").Value = 29

$wb.Save()
